$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Kelas" (class) column values for students
$ws.Range("O2").Value = "12 RPL A"
$ws.Range("O3").Value = "12 Mekatronika A"
$ws.Range("O4").Value = "12 Mekatronika A"
$ws.Range("O5").Value = "12 RPL A"
